$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ALNS schedule data for the Solution sheet, rows 2-10 (Staff_1..Staff_9), columns B-AC (Day 1..Day 28)
$scheduleData = @(
    @("DO","A1","A1","M3","M1","M3","A1","DO","A1","A1","A1","A1","M3","M3","M1","A1","A1","M3","DO","M1","M3","DO","A1","M3","M3","A1","A1","A1"),
    @("M3","DO","M1","M3","M2","M3","M3","M3","DO","M1","M3","M3","M2","M3","DO","M2","M2","M2","M3","M1","M3","M3","DO","M3","A1","M3","M3","A2"),
    @("DO","M3","A1","A1","A1","A1","M3","A1","M1","DO","M3","M3","A1","A1","A1","DO","M3","A1","A1","A1","M3","DO","M3","A1","M1","A1","A1","M3"),
    @("A1","DO","M3","M2","M3","M3","M3","DO","M3","M3","M3","M3","M1","M3","DO","M3","M1","M2","M2","M3","A2","A1","DO","M3","M3","M3","M3","M3"),
    @("A2","M1","M2","M2","M2","M1","DO","A1","M2","DO","M2","M1","M2","A1","A1","M1","DO","M2","M2","M1","A1","M3","M1","M3","A2","M2","M3","DO"),
    @("M3","A1","DO","A1","A1","A1","A1","M3","A1","A1","A1","A1","A1","DO","M3","A1","A1","A1","A1","A1","DO","M3","A1","A1","M1","A1","A1","DO"),
    @("A1","A2","DO","A2","M3","M2","M3","A1","DO","A1","A2","A1","M3","A1","A2","DO","A1","A1","M3","M2","A1","A2","A1","A2","A2","A2","M3","DO"),
    @("M1","M1","M3","DO","A2","A2","A2","M3","M3","M3","DO","M1","A2","M3","M3","M2","DO","M1","M3","A2","M3","M3","M2","M1","DO","M3","M3","A1"),
    @("M3","M3","A1","M3","M3","M3","DO","M3","A2","M1","M1","PH","M3","DO","M3","A1","M3","PH","A2","M2","DO","M3","M3","PH","M3","DO","M3","M1")
)

$startRow = 2
$startCol = 2  # Column B
for ($r = 0; $r -lt $scheduleData.Length; $r++) {
    $rowValues = $scheduleData[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $startCol + $c).Value = $rowValues[$c]
    }
}
